$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# LED forward-voltage variation means fewer of the 2K chip resistors are
# needed on this board revision: R86, R88 and R91 are dropped from the
# designator list, and the quantity for that BOM row drops from 7 to 4.
# The leading apostrophe preserves the cell's existing "stored as text"
# (quote-prefix) formatting instead of letting Excel re-infer a plain style.
$ws.Range("C6").Value = "'R85, R87, R89, R90"
$ws.Range("F6").Value = 4

# Leave the active selection on F7, matching the saved workbook state.
$ws.Range("F7").Select()
